$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8846
$ws1.Range("F3").Value = 96
$ws1.Range("F4").Value = 237
$ws1.Range("F5").Value = 102
$ws1.Range("F6").Value = 1464
$ws1.Range("F7").Value = 1401
$ws1.Range("F8").Value = 245
$ws1.Range("F9").Value = 46
$ws1.Range("F10").Value = 306
$ws1.Range("F11").Value = 89

# Sheet "全部类型" (fourth sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8847
$ws4.Range("F3").Value = 96
$ws4.Range("F4").Value = 237
$ws4.Range("F5").Value = 102
$ws4.Range("F6").Value = 1464
$ws4.Range("F7").Value = 1401
$ws4.Range("F8").Value = 245
$ws4.Range("F10").Value = 46
$ws4.Range("F11").Value = 306
$ws4.Range("F12").Value = 89
